$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("C1").Value = -0.69597319999999996
$ws.Range("C2").Value = 5.0163497000000001
$ws.Range("C3").Value = 0.17579020000000001

$ws.Range("C3").Select()
